# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    3  = 5003
    4  = 4
    5  = 7267
    6  = 44
    7  = 66
    9  = 595
    10 = 64
    11 = 16
    12 = 4248
    13 = 1708
    15 = 93
    16 = 2849
    18 = 560
    19 = 195
    20 = 452
    21 = 413
    22 = 438
    23 = 273
    24 = 79
    26 = 1142
    27 = 83
    28 = 1332
    29 = 100
    30 = 568
    31 = 18
    32 = 509
    34 = 48
    36 = 2660
    37 = 680
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
